$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: update the Admin password value (B2) ---
$ws1.Range("B2").Value = "testadmin123"

# --- Add Sheet2 after Sheet1 with new datadriven test data ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "SearchData"
$ws2.Range("A2").Value = "Bag"
$ws2.Range("A3").Value = "abab"
$ws2.Range("A5").Value = "belt"
$ws2.Range("A4").Value = "Jacket"

$ws2.Columns.Item(1).ColumnWidth = 22.5

# --- Selections: Sheet1 -> B2, Sheet2 (active) -> A4 ---
$ws1.Range("B2").Select()
$ws2.Range("A4").Select()
